# Delete row 435 (the post "「あなたの声を聞かせてください。私の心を彩ってください」")
# which removes that entire row and shifts all subsequent rows up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows.Item(435).Delete()
